$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value2 = $val
    $cell.Style = $origStyle
}

$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '60.407.72', '  +2.44%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '3.032.22', '  +1.12%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.00', '  -0.16%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '569.02', '  +1.27%  '),
    @(6, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '141.17', '  +2.23%  '),
    @(7, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.00', '  -0.08%  '),
    @(8, 'LidoStakedEther', 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth', '3.035.04', '  +1.40%  '),
    @(9, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.523', '  +0.98%  '),
    @(10, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.137', '  +3.31%  '),
    @(11, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '5.41', '  +11.02%  '),
    @(12, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.457', '  -0.19%  '),
    @(13, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0000235', '  +2.04%  '),
    @(14, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '34.50', '  +2.06%  '),
    @(15, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.123', '  -0.26%  '),
    @(16, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '3.531.12', '  +1.09%  '),
    @(17, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '7.19', '  +2.54%  '),
    @(18, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '3.027.58', '  +1.10%  '),
    @(19, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '60.394.50', '  +2.41%  '),
    @(20, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '446.53', '  +4.38%  '),
    @(21, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '13.79', '  +1.12%  '),
    @(22, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.731', '  +2.33%  '),
    @(23, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '7.21', '  +0.61%  '),
    @(24, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '13.54', '  +0.23%  '),
    @(25, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '81.21', '  +0.85%  '),
    @(26, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.00', '  +0.18%  '),
    @(27, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '2.28', '  +8.13%  '),
    @(28, 'FirstDigitalUSD', 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd', '1.00', '  -0.18%  '),
    @(29, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '8.02', '  +3.22%  '),
    @(30, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '2.58', '  +1.73%  '),
    @(31, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '6.39', '  +3.95%  '),
    @(32, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '26.28', '  +1.83%  '),
    @(33, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.106', '  +7.70%  '),
    @(34, 'PEPE', 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe', '0.0₃0807', '  +6.46%  '),
    @(35, 'Mantle', 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt', '1.02', '  +4.85%  '),
    @(36, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '6.02', '  +4.37%  '),
    @(37, 'Stacks', 'https://coinranking.com/coin/mMPrMcB7+stacks-stx', '2.13', '  +1.47%  '),
    @(38, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '49.44', '  +1.14%  '),
    @(39, 'dogwifhat', 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif', '2.94', '  +7.00%  '),
    @(40, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '8.74', '  -1.38%  '),
    @(41, 'Bittensor', 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao', '407.09', '  +2.98%  '),
    @(42, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.0357', '  +1.85%  '),
    @(43, 'Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '2.766.98', '  +1.36%  '),
    @(44, 'Kaspa', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', '0.107', '  -1.13%  '),
    @(45, 'TheGraph', 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt', '0.259', '  +4.74%  '),
    @(46, 'Arweave', 'https://coinranking.com/coin/7XWg41D1+arweave-ar', '37.01', '  +14.57%  '),
    @(47, 'USDe', 'https://coinranking.com/coin/exbfr2U-0+usde-usde', '0.999', '  -0.01%  '),
    @(48, 'Fetch.AI', 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet', '2.07', '  +1.56%  '),
    @(49, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '123.30', '  -1.80%  '),
    @(50, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.110', '  +0.54%  '),
    @(51, 'InjectiveProtocol', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj', '23.66', '  +0.87%  '),
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 2).Value2 = $entry[1]
    $ws.Cells.Item($row, 3).Value2 = $entry[2]
    Set-TextValue $row 4 $entry[3]
    $ws.Cells.Item($row, 5).Value2 = $entry[4]
}

Write-Output "Updated cryptos list rows 2-51"